$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Color info")

# Delete the "R: nnn, G: nnn, B:nnn" / "#ffffff" row (row 5), shifting the
# "Spectral distribution" row (old row 6) up into its place.
$ws.Rows.Item(5).Delete()
